# Swap the data for rows 4 and 5 in columns A, B, D, E, F, G, H, I
# (columns J onward already hold identical values in both rows, so no change needed there)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "I")

foreach ($col in $cols) {
    $addr4 = $col + "4"
    $addr5 = $col + "5"
    $cell4 = $ws.Range($addr4)
    $cell5 = $ws.Range($addr5)

    $val4 = $cell4.Value2
    $val5 = $cell5.Value2

    $cell4.Value2 = $val5
    $cell5.Value2 = $val4
}
